{"js": "// Apply the text replacements described by the diff: the date line and\n// the 25 two-digit multiplication problems in the table.\nconst replacements = [\n  [\"2025-05-02 Friday\", \"2025-05-03 Saturday\"],\n  [\"26\u00d761=\", \"71\u00d712=\"],\n  [\"25\u00d788=\", \"39\u00d736=\"],\n  [\"35\u00d742=\", \"55\u00d715=\"],\n  [\"32\u00d713=\", \"20\u00d753=\"],\n  [\"30\u00d764=\", \"14\u00d718=\"],\n  [\"89\u00d764=\", \"56\u00d760=\"],\n  [\"51\u00d760=\", \"79\u00d711=\"],\n  [\"77\u00d712=\", \"56\u00d791=\"],\n  [\"38\u00d734=\", \"33\u00d726=\"],\n  [\"41\u00d722=\", \"76\u00d749=\"],\n  [\"47\u00d778=\", \"68\u00d750=\"],\n  [\"84\u00d749=\", \"30\u00d750=\"],\n  [\"37\u00d789=\", \"28\u00d769=\"],\n  [\"82\u00d754=\", \"97\u00d754=\"],\n  [\"81\u00d772=\", \"14\u00d747=\"],\n  [\"55\u00d749=\", \"84\u00d715=\"],\n  [\"14\u00d721=\", \"42\u00d771=\"],\n  [\"80\u00d727=\", \"93\u00d723=\"],\n  [\"72\u00d751=\", \"46\u00d796=\"],\n  [\"27\u00d759=\", \"51\u00d743=\"],\n  [\"21\u00d724=\", \"16\u00d791=\"],\n  [\"28\u00d792=\", \"99\u00d770=\"],\n  [\"83\u00d748=\", \"53\u00d757=\"],\n  [\"44\u00d751=\", \"90\u00d758=\"],\n  [\"15\u00d741=\", \"91\u00d720=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the text replacements described by the diff: the date line and\n# the 25 two-digit multiplication problems in the table.\n$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-05-02 Friday\", \"2025-05-03 Saturday\"),\n  @(\"26\u00d761=\", \"71\u00d712=\"),\n  @(\"25\u00d788=\", \"39\u00d736=\"),\n  @(\"35\u00d742=\", \"55\u00d715=\"),\n  @(\"32\u00d713=\", \"20\u00d753=\"),\n  @(\"30\u00d764=\", \"14\u00d718=\"),\n  @(\"89\u00d764=\", \"56\u00d760=\"),\n  @(\"51\u00d760=\", \"79\u00d711=\"),\n  @(\"77\u00d712=\", \"56\u00d791=\"),\n  @(\"38\u00d734=\", \"33\u00d726=\"),\n  @(\"41\u00d722=\", \"76\u00d749=\"),\n  @(\"47\u00d778=\", \"68\u00d750=\"),\n  @(\"84\u00d749=\", \"30\u00d750=\"),\n  @(\"37\u00d789=\", \"28\u00d769=\"),\n  @(\"82\u00d754=\", \"97\u00d754=\"),\n  @(\"81\u00d772=\", \"14\u00d747=\"),\n  @(\"55\u00d749=\", \"84\u00d715=\"),\n  @(\"14\u00d721=\", \"42\u00d771=\"),\n  @(\"80\u00d727=\", \"93\u00d723=\"),\n  @(\"72\u00d751=\", \"46\u00d796=\"),\n  @(\"27\u00d759=\", \"51\u00d743=\"),\n  @(\"21\u00d724=\", \"16\u00d791=\"),\n  @(\"28\u00d792=\", \"99\u00d770=\"),\n  @(\"83\u00d748=\", \"53\u00d757=\"),\n  @(\"44\u00d751=\", \"90\u00d758=\"),\n  @(\"15\u00d741=\", \"91\u00d720=\")\n)\n\nforeach ($pair in $replacements) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
